$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row to append: row 90
# A90: Date (text, same shape as existing date cells)
# B90: Game name
# C90: ModCount (numeric)

# Force column A to be treated as text so "2026/02/08" isn't
# auto-converted into a date serial number, then restore the
# cell's style to "Normal" (clears the temporary text number
# format) before re-applying the same center/center alignment
# used by the rest of the data rows.
$ws.Range("A90").NumberFormat = "@"
$ws.Range("A90").Value = "2026/02/08"
$ws.Range("A90").Style = "Normal"
$ws.Range("A90").HorizontalAlignment = -4108
$ws.Range("A90").VerticalAlignment = -4108

$ws.Range("B90").Value = "逃离鸭科夫"
$ws.Range("B90").HorizontalAlignment = -4108
$ws.Range("B90").VerticalAlignment = -4108

$ws.Range("C90").Value = 1184
$ws.Range("C90").HorizontalAlignment = -4108
$ws.Range("C90").VerticalAlignment = -4108
